$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Add new row 8 with version history entry
$ws.Range("A8").Value = Get-Date -Year 2025 -Month 3 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("A8").NumberFormat = "m/d/yy"

$ws.Range("B8").Value = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 16 -Minute 19 -Second 0
$ws.Range("B8").NumberFormat = "h:mm:ss"

$ws.Range("C8").Value = "Futconnect2703 1619"
$ws.Range("D8").Value = "Politicas em mensalidades, jogos e sócios - feito o commit"

# Update selection to reflect the next empty cell after the new row
$ws.Range("D9").Select()
